# This script updates column B ("alpha") values in Sheet1 to reflect
# the new alpha parameter added to the LDA run (commit: "adding alpha to
# lda(), adding alphas.xlsx spreadsheet"). Only specific rows change value;
# all other cells remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 3
$ws.Range("B8").Value = 4
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 2
$ws.Range("B14").Value = 3
$ws.Range("B17").Value = 3
$ws.Range("B19").Value = 5
$ws.Range("B20").Value = 2
$ws.Range("B21").Value = 2
$ws.Range("B24").Value = 2
$ws.Range("B25").Value = 3
$ws.Range("B26").Value = 2
$ws.Range("B28").Value = 3
$ws.Range("B29").Value = 2
$ws.Range("B30").Value = 4
$ws.Range("B31").Value = 2
$ws.Range("B35").Value = 5
$ws.Range("B36").Value = 4
$ws.Range("B37").Value = 5
$ws.Range("B38").Value = 5
$ws.Range("B39").Value = 3
$ws.Range("B40").Value = 2
$ws.Range("B41").Value = 3
$ws.Range("B42").Value = 5
$ws.Range("B43").Value = 5
$ws.Range("B44").Value = 4
$ws.Range("B46").Value = 4
$ws.Range("B49").Value = 2
$ws.Range("B52").Value = 4
$ws.Range("B53").Value = 3
$ws.Range("B56").Value = 5
$ws.Range("B57").Value = 5
$ws.Range("B58").Value = 5
$ws.Range("B61").Value = 5
$ws.Range("B62").Value = 5
$ws.Range("B63").Value = 4
$ws.Range("B64").Value = 2
$ws.Range("B65").Value = 5
$ws.Range("B66").Value = 5
$ws.Range("B67").Value = 5
$ws.Range("B68").Value = 5
$ws.Range("B69").Value = 2
$ws.Range("B70").Value = 5
$ws.Range("B71").Value = 5
$ws.Range("B73").Value = 3
$ws.Range("B74").Value = 3
$ws.Range("B76").Value = 4
$ws.Range("B77").Value = 5
$ws.Range("B78").Value = 4
$ws.Range("B79").Value = 5
$ws.Range("B80").Value = 4
$ws.Range("B81").Value = 3
$ws.Range("B82").Value = 2
$ws.Range("B83").Value = 5
$ws.Range("B84").Value = 4
$ws.Range("B85").Value = 1
$ws.Range("B87").Value = 5
$ws.Range("B88").Value = 1
$ws.Range("B89").Value = 4
$ws.Range("B90").Value = 1
$ws.Range("B91").Value = 1
$ws.Range("B92").Value = 1
$ws.Range("B95").Value = 3
$ws.Range("B97").Value = 2
$ws.Range("B98").Value = 3
$ws.Range("B99").Value = 2
$ws.Range("B100").Value = 5
$ws.Range("B101").Value = 4
